# Apply AkWarm/NEEP Library update (9/27/2018 AkWarm library) to the Misc Info sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LibVersion date: 4/4/2018 -> 9/27/2018 (serial 43194 -> 43370)
$ws.Range("B2").Value = 43370

# RegSurcharge: 0.0032 -> 0.004 (stored as float32-precision doubles)
$ws.Range("D2").Value = 0.00400000018998981

# RegSurchargeElectric: 0.0009 -> 0.000978 (stored as float32-precision doubles)
$ws.Range("E2").Value = 0.000977999996393919

# Update MiscNotes text: "July 2012" -> "July 2018" and gas RCC ".32%" -> ".40%"
$ws.Range("H2").Value = "Inflation factors and discount rate from 2011 FEMP;  Regulatory surcharge is now correctly modeled as a % for gas utilities and a `$/kWh surcharge for electric utilities.  PCE 100% funding in effect July 2018. Jan 2018 gas RCC is .40% of total bill - fuel costs use price before taxes added as program now does the math"

Write-Output "Updated B2, D2, E2, H2 on Misc_Info sheet"
